# Reported two window resize issues related to Transaction UI
# This script reproduces the OOXML diff: two new "Issue" rows describing
# Transaction UI resize bugs are appended to Sheet1's issue-tracker table,
# the "#" counter column (A) is back-filled for the previously-unnumbered
# rows 33-41, and a new run of counter-only rows (44-53) is added below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# ---------------------------------------------------------------------
# 1. Back-fill column A ("#") for existing rows 33-41, which previously
#    had no value in column A (spans was "2:9").
# ---------------------------------------------------------------------
$ws.Range("A33").Value = 31
$ws.Range("A34").Value = 32
$ws.Range("A35").Value = 33
$ws.Range("A36").Value = 34
$ws.Range("A37").Value = 35
$ws.Range("A38").Value = 36
$ws.Range("A39").Value = 37
$ws.Range("A40").Value = 38
$ws.Range("A41").Value = 39

# Row 41's wrapped content grows taller once considered together with the
# two new rows below it; its custom height changes from 67.5 to 94.5.
$ws.Rows.Item(41).RowHeight = 94.5

# ---------------------------------------------------------------------
# 2. Add the two new "Transaction UI" issue rows (42 and 43).
#    The shared-string table records new unique strings in first-seen
#    order, so cells are populated in the same column-major order the
#    original authoring tool used: B42/B43 ("Transaction UI"), then the
#    "Observed Result" column (C42, C43), then the "Additional remark"
#    column (D42, D43), to land on the exact shared-string indices.
# ---------------------------------------------------------------------
$ws.Range("A42").Value = 40
$ws.Range("B42").Value = "Transaction UI"

$ws.Range("A43").Value = 41
$ws.Range("B43").Value = "Transaction UI"

$ws.Range("C42").Value = "Resizing main window does not resize Transaction panels accordingly."
$ws.Range("C43").Value = "Clicking View button, resets size of main window."

$ws.Range("D42").Value = "Case 1:`nSelect Reports->Transactions`nMaximize main app window`nObserve that Transaction panel does not resizeaccordingly`nCase 2:`nInstead of maximizing the window, reduce the size.`nObserve that the UI is clipped."
$ws.Range("D43").Value = "Maximize main app window`nSelect Reports->Transactions`nClick search`nSelect transaction from list and click View button`nObserve that the main window size is reset"

$ws.Range("F42").Value = "Nikhil"
$ws.Range("F43").Value = "Nikhil"

$ws.Rows.Item(42).RowHeight = 135
$ws.Rows.Item(43).RowHeight = 90

# ---------------------------------------------------------------------
# 3. Add the trailing counter-only rows 44-53 (column A keeps counting
#    upward; no other columns are populated on these rows).
# ---------------------------------------------------------------------
$ws.Range("A44").Value = 42
$ws.Range("A45").Value = 43
$ws.Range("A46").Value = 44
$ws.Range("A47").Value = 45
$ws.Range("A48").Value = 46
$ws.Range("A49").Value = 47
$ws.Range("A50").Value = 48
$ws.Range("A51").Value = 49
$ws.Range("A52").Value = 50
$ws.Range("A53").Value = 51

# ---------------------------------------------------------------------
# 4. Update the view: the sheet had scrolled down to show the new rows
#    and the last touched cell was I42.
# ---------------------------------------------------------------------
$ws.Range("I42").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
